$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B/C/D columns for rows 2-27 with new ticker values
$ws.Range("B2").Value = "NSE:BHARATGEAR"
$ws.Range("C2").Value = "NSE:3MINDIA"
$ws.Range("D2").Value = "NSE:ASIANPAINT"
$ws.Range("B3").Value = "NSE:DHANUKA"
$ws.Range("C3").Value = "NSE:ASAHIINDIA"
$ws.Range("D3").Value = "NSE:JSWSTEEL"
$ws.Range("B4").Value = "NSE:EBBETF0431"
$ws.Range("C4").Value = "NSE:AWL"
$ws.Range("B5").Value = "NSE:GLAND"
$ws.Range("C5").Value = "NSE:BAJAJFINSV"
$ws.Range("B6").Value = "NSE:HIRECT"
$ws.Range("C6").Value = "NSE:BANDHANBNK"
$ws.Range("B7").Value = "NSE:LUPIN"
$ws.Range("C7").Value = "NSE:BEARDSELL"
$ws.Range("B8").Value = "NSE:MEDPLUS"
$ws.Range("C8").Value = "NSE:BORORENEW"
$ws.Range("B9").Value = "NSE:PUNJABCHEM"
$ws.Range("C9").Value = "NSE:DALBHARAT"
$ws.Range("B10").Value = "NSE:RAMRAT"
$ws.Range("C10").Value = "NSE:DELTACORP"
$ws.Range("C11").Value = "NSE:ESABINDIA"
$ws.Range("C12").Value = "NSE:GANGAFORGE"
$ws.Range("C13").Value = "NSE:GTPL"
$ws.Range("C14").Value = "NSE:KAMATHOTEL"
$ws.Range("C15").Value = "NSE:KFINTECH"
$ws.Range("C16").Value = "NSE:MAHLIFE"
$ws.Range("C17").Value = "NSE:MANAKCOAT"
$ws.Range("C18").Value = "NSE:MASFIN"
$ws.Range("C19").Value = "NSE:MOIL"
$ws.Range("C20").Value = "NSE:NESTLEIND"
$ws.Range("C21").Value = "NSE:NGIL"
$ws.Range("C22").Value = "NSE:NUVAMA"
$ws.Range("C23").Value = "NSE:NV20BEES"
$ws.Range("C24").Value = "NSE:OLECTRA"
$ws.Range("C25").Value = "NSE:PDSL"
$ws.Range("C26").Value = "NSE:RAINBOW"
$ws.Range("C27").Value = "NSE:ROSSARI"

# Clear cells that previously had values in E/F columns but are now blank
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("E10").ClearContents()

# Remove rows 28-42 (workbook shrinks from A1:F42 to A1:F27)
$ws.Range("A28:F42").EntireRow.Delete()

Write-Host "edit complete"
